$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.708.03"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.299.30"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'301.09"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").Value = "'96.02"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("D7").Value = "'0.510"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "'34.73"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").Value = "'19.28"
$ws.Range("E11").Value = "  +5.52%  "
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'6.81"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "2.651.37"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "2.294.82"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "'0.784"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "42.644.60"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "'12.29"
$ws.Range("E19").Value = "  -6.12%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").Value = "'67.67"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").Value = "'2.40"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").Value = "'24.60"
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("D28").Value = "'164.85"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").Value = "'32.09"
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'4.97"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").Value = "'4.43"
$ws.Range("E35").Value = "  -8.21%  "
$ws.Range("D36").Value = "'0.0699"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("D38").Value = "'0.0999"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").Value = "'2.71"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").Value = "'19.80"
$ws.Range("E42").Value = "  +9.53%  "
$ws.Range("D43").Value = "1.961.89"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").Value = "'10.46"
$ws.Range("E44").Value = "  +4.66%  "
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("E46").Value = "  -3.18%  "
$ws.Range("D47").Value = "'2.75"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").Value = "2.524.66"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "'53.16"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").Value = "'71.63"
$ws.Range("E51").Value = "  -0.46%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'235.20"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'2.24"
$ws.Range("E24").Value = "  +3.10%  "
